# Regenerate save_data to use K (column G) instead of Strike#.
# New K values (calculated/regenerated) are written into column G for rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 0
    4  = 2
    5  = 0
    6  = 3
    7  = 1
    8  = 0
    9  = 2
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 3
    16 = 0
    17 = 1
    18 = 2
    19 = 0
    20 = 1
    21 = 1
    22 = 0
    23 = 0
    24 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
